$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.539.84'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '1.803.07'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.603'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.96'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +16.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.292'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0995'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.07%  '
$ws.Range("D12").Value = '2.062.29'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.812.47'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '34.523.36'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '0.0₃0768'
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.02%  '
$ws.Range("E24").Value = '  -2.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.121'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '87.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.77%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = '1.316.43'
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0187'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.37%  '
$ws.Range("E42").Value = '  +4.91%  '
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.938'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0518'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.73%  '
$ws.Range("D47").Value = '1.963.48'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("E51").Value = '  +0.44%  '
